# Applies the "adding Detainees n, some formatting edits" commit to the
# CCJ_quantified_values workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Clear the placeholder "-" text out of the Mean/Median/Max/Min cells
#    that only ever held a dash (they stay numeric-formatted/empty now).
# ---------------------------------------------------------------------
$dashCells = @(
    "I7",
    "H8", "J8", "K8",
    "I9", "J9",
    "J10", "K10",
    "H11", "J11", "K11",
    "H12", "I12", "J12", "K12",
    "H13", "J13", "K13",
    "I14", "J14"
)
foreach ($cellAddr in $dashCells) {
    $ws.Range($cellAddr).ClearContents()
}

# ---------------------------------------------------------------------
# 2) Fill in the "impact" column (M) that was previously left blank for
#    several rows.
# ---------------------------------------------------------------------
$ws.Range("M8").Value = "negative"
$ws.Range("M9").Value = "negative"
$ws.Range("M10").Value = "negative"
$ws.Range("M11").Value = "negative"

$ws.Range("M15").Value = "cost"
$ws.Range("M16").Value = "cost"
$ws.Range("M17").Value = "cost"
$ws.Range("M18").Value = "cost"

# ---------------------------------------------------------------------
# 3) Add the missing source / link for the "Total Funding for CCJ" row.
# ---------------------------------------------------------------------
$ws.Range("O18").Value = "CCJ Dashboard"
$ws.Range("P18").Value = "https://cook-dashboard.loyolaccj.org/jail/admissions?utm_source=chatgpt.com"

# ---------------------------------------------------------------------
# 4) New "Weight" rows for n Detainees / n Society.
# ---------------------------------------------------------------------
$ws.Range("A19").Value = "Weight"
$ws.Range("B19").Value = "n Detainees"
$ws.Range("D19").Value = "people"
$ws.Range("G19").Value = 33945

$ws.Range("A20").Value = "Weight"
$ws.Range("B20").Value = "n Society "
$ws.Range("D20").Value = "people"

# Highlight the "n Society" subcomponent cell in yellow (new fill added
# to the workbook's style table).
$ws.Range("B20").Interior.Color = 65535

# ---------------------------------------------------------------------
# 5) Bold the header row (row 1).
# ---------------------------------------------------------------------
$ws.Range("A1:R1").Font.Bold = $true
